# Swap the `name` given to the Pearson/BTEC logo inline pictures that
# live in this document's headers and footers:
#
#   footers (Pearson logo, descr="...PearsonLogo.png"): image1.png -> image2.png
#   headers (BTEC logo,    descr="BTec_Logo-Orange"):   image2.jpg -> image1.jpg
#
# Word's InlineShape object doesn't track a picture's on-disk media
# filename anywhere obvious, so we key off the existing AlternativeText
# (== the <wp:docPr descr="..."/> value, which is untouched by this
# change) together with the shape's current Name to decide what the
# new Name should become.

$d = $word.ActiveDocument

function Rename-LogoShape($headerFooter) {
    if (-not $headerFooter.Exists) { return }
    $shapes = $headerFooter.Range.InlineShapes
    for ($shapeIdx = 1; $shapeIdx -le $shapes.Count; $shapeIdx++) {
        # Re-fetch through Selection each time -- renaming works
        # reliably for both header- and footer-hosted pictures when
        # driven off $word.Selection.InlineShapes.
        $shapes.Item($shapeIdx).Range.Select()
        $shape = $word.Selection.InlineShapes.Item(1)

        $descr = $shape.AlternativeText

        if ($descr -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shape.Name = "image2.png"
        } elseif ($descr -eq "BTec_Logo-Orange") {
            $shape.Name = "image1.jpg"
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($hfIdx = 1; $hfIdx -le 3; $hfIdx++) {
        Rename-LogoShape $sec.Headers.Item($hfIdx)
        Rename-LogoShape $sec.Footers.Item($hfIdx)
    }
}
